$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (the "0" placeholder cell with special border/bold/center style);
# this shifts the shared-string cell up from A2 to A1 and drops the row-1 style.
$ws.Rows.Item(1).Delete()

# Update the remaining cell's text to the reformatted (pretty-printed) questions payload.
$ws.Range("A1").Value = @'
questions = [
    {
        "title": "You have just closed an opportunity in Salesforce CRM and want to set up a follow-up activity for a few weeks in the future that is properly connected to the opportunity record.  How should you do this?",
        "ques_type": 2,
        "options": [
            "Enter the follow-up activity information into the opportunity record manually.",
            "Log the follow-up task and set a due date.",
            "Create a new custom object for follow-up activities and manually link it to the opportunity record.",
            "Set a call reminder on an external calendar."
        ],
        "score": "Log the follow-up task and set a due date."
    },
    {
        "title": "You manage accounts and contacts for a sales team that has identified a new account, a medium-sized business, that it needs to create contacts for. Which steps should you take?",
        "ques_type": 15,
        "options": [
            "Create a new account.",
            "Assign a primary contact for the new account.",
            "Add new contacts to the account with accurate contact information.",
            "Track all interactions with each new contact.",
            "Maintain a detailed record of the new account's history.",
            "Merge the new account with an existing account in Salesforce."
        ],
        "score": [
            "Create a new account.",
            "Assign a primary contact for the new account.",
            "Add new contacts to the account with accurate contact information."
        ]
    },
    {
        "title": "You have been working with a lead for several weeks, and after numerous calls and meetings, they have agreed to purchase your product. You need to update their information in Salesforce CRM. Which is the first step you should take?",
        "ques_type": 2,
        "options": [
            "Create a new Opportunity.",
            "Convert the Lead to an Opportunity.",
            "Add the Opportunity to the existing account.",
            "Update the lead status to Closed Won."
        ],
        "score": "Convert the Lead to an Opportunity."
    },
    {
        "title": "Your sales team is having trouble identifying decision-makers within their accounts, which is hindering their ability to close deals. As the Salesforce administrator, you need to find a solution to help them identify decision-makers more easily. Which action should you take?",
        "ques_type": 2,
        "options": [
            "Add contact roles to each account.",
            "Create custom fields on the account object to track decision-maker information.",
            "Use the Account Hierarchy feature to identify decision-makers.",
            "Train your sales team on better prospecting techniques and strategies for identifying decision-makers."
        ],
        "score": "Add contact roles to each account."
    }
]
'@

# Re-fit the row height back to default (setting a long multi-line value auto-expands it).
$ws.Rows.Item(1).AutoFit()
